# "actualizar boton categoria seleccion"
# - Remove the now-unused "TRABAJO" category row from Categorias.
# - Bump Contador for several Palabras rows (selection counters).
# - Append two newly-selected words (guineo, coco) to Palabras.

$wb = $excel.ActiveWorkbook

# --- Sheet: Categorias -------------------------------------------------
$wsCat = $wb.Worksheets.Item("Categorias")
# Row 6 = ID_Categoria 5 / TRABAJO / 2025-03-06 18:24:32
$wsCat.Rows.Item(6).Delete()

# --- Sheet: Palabras -----------------------------------------------------
$wsPal = $wb.Worksheets.Item("Palabras")

# Increment the "Contador" (column D) for the words that were selected again.
$wsPal.Cells.Item(2, 4).Value = 4    # pollo:    3 -> 4
$wsPal.Cells.Item(3, 4).Value = 5    # frutas:   4 -> 5
$wsPal.Cells.Item(4, 4).Value = 7    # toalla:   6 -> 7
$wsPal.Cells.Item(6, 4).Value = 1    # aceite:   0 -> 1
$wsPal.Cells.Item(8, 4).Value = 10   # desayuno: 9 -> 10
$wsPal.Cells.Item(11, 4).Value = 1   # caliente: 0 -> 1
$wsPal.Cells.Item(16, 4).Value = 1   # cena:     0 -> 1
$wsPal.Cells.Item(22, 4).Value = 1   # jugo:     0 -> 1
$wsPal.Cells.Item(23, 4).Value = 1   # dulce:    0 -> 1

# Append the two new words picked in this category.
$wsPal.Cells.Item(25, 1).Value = 24
$wsPal.Cells.Item(25, 2).Value = 1
$wsPal.Cells.Item(25, 3).Value = "guineo"
$wsPal.Cells.Item(25, 4).Value = 0
$wsPal.Cells.Item(25, 5).Value = "2025-03-06 18:42:07"

$wsPal.Cells.Item(26, 1).Value = 25
$wsPal.Cells.Item(26, 2).Value = 1
$wsPal.Cells.Item(26, 3).Value = "coco"
$wsPal.Cells.Item(26, 4).Value = 1
$wsPal.Cells.Item(26, 5).Value = "2025-03-06 18:53:22"
